$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: fill in the new journal entry (12.03.2020, 16:00 -> 16:15, GitHub / Bataille Navale / Planifications des sprint / Création des sprints pour chaques semaines)
$ws.Range("B8").Value = 43902
$ws.Range("C8").Value = 0.66666666666666663
$ws.Range("D8").Value = 0.67708333333333337
$ws.Range("F8").Value = "GitHub"
$ws.Range("G8").Value = "Bataille Navale"
$ws.Range("H8").Value = "Planifications des sprint"
$ws.Range("I8").Value = "Création des sprints pour chaques semaines"

# Update the selected cell/active cell in the sheet view to D9
$ws.Range("D9").Select() | Out-Null
